# commit and push by Abhishek k- TC_15
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the "Vendor" sheet to "Vendors"
# ---------------------------------------------------------------------
$wsVendors = $wb.Worksheets.Item("Vendor")
$wsVendors.Name = "Vendors"

# ---------------------------------------------------------------------
# 2) Organizations sheet: fix the Org_Name value for TC_03
# ---------------------------------------------------------------------
$wsOrg = $wb.Worksheets.Item("Organizations")
$wsOrg.Activate()
$wsOrg.Range("C2").Value = "IAmFromBidar"
$wsOrg.Range("C2").Select()

# ---------------------------------------------------------------------
# 3) Vendors sheet: add the new TC_12 test-case block (fill the ID down
#    first, then the rest of the header row)
# ---------------------------------------------------------------------
$wsVendors.Range("A7").Value = "TC_ID"
$wsVendors.Range("B7").Value = "TC_Name"
$wsVendors.Range("A8").Value = "TC_12"
$wsVendors.Range("C7").Value = "Vendors_Name"
$wsVendors.Range("D7").Value = "Product_Name"
$wsVendors.Activate()
$wsVendors.Range("D7").Select()

# ---------------------------------------------------------------------
# 4) Products sheet: update vendor/product values for TC_05, widen col C,
#    and leave this sheet as the active one
# ---------------------------------------------------------------------
$wsProducts = $wb.Worksheets.Item("Products")
$wsProducts.Activate()
$wsProducts.Range("D2").Value = "Oneplus_11r"
$wsProducts.Range("C2").Value = "Abhishek Kelusker"
$wsProducts.Columns.Item(3).ColumnWidth = 17.14
$wsProducts.Range("D4").Select()

# ---------------------------------------------------------------------
# 5) Leave "Products" as the active tab/sheet when the workbook is saved
# ---------------------------------------------------------------------
$wsProducts.Activate()
